# Add test kohde without any osapuoli
# Adds a new "rakennus" (building) row to the "R1 rakennus" sheet and a
# matching "osoite" (address) row to the "R3 osoite" sheet, for a building
# that has no owner/occupant ("osapuoli") rows anywhere.

$wb = $excel.ActiveWorkbook

$wsRakennus = $wb.Worksheets.Item("R1 rakennus")
$wsOsoite   = $wb.Worksheets.Item("R3 osoite")

# --- R1 rakennus: new row 7 --------------------------------------------
$wsRakennus.Range("A7").Value = "145678901C"
$wsRakennus.Range("B7").Value = "398"
$wsRakennus.Range("C7").Value = 20
$wsRakennus.Range("D7").Value = "39800300010001"
$wsRakennus.Range("E7").Value = "Salattu"
$wsRakennus.Range("F7").Value = "39800300010001"
$wsRakennus.Range("G7").Value = 1
$wsRakennus.Range("H7").Value = 1
$wsRakennus.Range("J7").Value = "15230"
$wsRakennus.Range("K7").Value = 20210101
$wsRakennus.Range("L7").Value = 1
$wsRakennus.Range("N7").Value = 260
$wsRakennus.Range("O7").Value = 2
$wsRakennus.Range("S7").Value = "01"
$wsRakennus.Range("T7").Value = 20210101
$wsRakennus.Range("U7").Value = "011"
$wsRakennus.Range("V7").Value = 1
$wsRakennus.Range("W7").Value = 1
$wsRakennus.Range("X7").Value = 6762345
$wsRakennus.Range("Y7").Value = 431007
$wsRakennus.Range("Z7").Value = 5

# --- R3 osoite: new row 7 ------------------------------------------------
$wsOsoite.Range("A7").Value = "145678901C"
$wsOsoite.Range("B7").Value = "398"
$wsOsoite.Range("C7").Value = 1
$wsOsoite.Range("D7").Value = "Halmekatu"
$wsOsoite.Range("F7").Value = 99
$wsOsoite.Range("G7").Value = "15230"
$wsOsoite.Range("H7").Value = "LAHTI"
$wsOsoite.Range("I7").Value = "LAHTIS"

# --- Selection / active tab ----------------------------------------------
# The building sheet becomes the active tab, selection parked one row below
# the newly added data (mirrors Excel's behaviour after data entry).
$wsRakennus.Range("A8").Select()
$wsOsoite.Range("A8").Select()
$wsRakennus.Activate()
